{"js": "// \"add adss on user\"\n//\n// 1. Justify (\"both\" / justified) the \"Are you a family member or a\n//    worker of a farmer...\" paragraph.\n// 2. Justify (\"both\" / justified) the \"I hereby certify...\" paragraph,\n//    and merge its two runs (which were split around a stray\n//    `_GoBack` bookmark) back into a single contiguous run of text.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst familyMemberNeedle = \"family member or a worker of a farmer\";\nconst certifyNeedle = \"I hereby certify that the forgoing information\";\n\nlet familyMemberPara = null;\nlet certifyPara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  if (familyMemberPara === null && text.indexOf(familyMemberNeedle) !== -1) {\n    familyMemberPara = para;\n  }\n  if (certifyPara === null && text.indexOf(certifyNeedle) !== -1) {\n    certifyPara = para;\n  }\n}\n\n// 1) \"Are you a family member or a worker of a farmer ...\" -> justified.\nif (familyMemberPara) {\n  familyMemberPara.alignment = Word.Alignment.justified;\n}\n\n// 2) \"I hereby certify ...\" -> justified, and collapse the bookmark-split\n// runs into a single run by re-writing the paragraph's full text in place\n// (this removes the now-pointless `_GoBack` bookmark along the way).\nif (certifyPara) {\n  certifyPara.alignment = Word.Alignment.justified;\n  const fullText = certifyPara.text;\n  certifyPara.insertText(fullText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# \"add adss on user\"\n#\n# 1. Justify (\"both\" / justified) the \"Are you a family member or a\n#    worker of a farmer...\" paragraph.\n# 2. Justify (\"both\" / justified) the \"I hereby certify...\" paragraph,\n#    and merge its two runs (which were split around a stray\n#    `_GoBack` bookmark) back into a single contiguous run of text.\n\n$wdAlignParagraphJustify = 3\n\n$d = $word.ActiveDocument\n\n$familyMemberNeedle = \"family member or a worker of a farmer\"\n$certifyNeedle = \"I hereby certify that the forgoing information\"\n\n$familyMemberPara = $null\n$certifyPara = $null\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($null -eq $familyMemberPara -and $t -like \"*$familyMemberNeedle*\") {\n        $familyMemberPara = $p\n    }\n    if ($null -eq $certifyPara -and $t -like \"*$certifyNeedle*\") {\n        $certifyPara = $p\n    }\n    if ($familyMemberPara -and $certifyPara) { break }\n}\n\n# 1) \"Are you a family member or a worker of a farmer ...\" -> justified.\nif ($familyMemberPara) {\n    $familyMemberPara.Alignment = $wdAlignParagraphJustify\n}\n\n# 2) \"I hereby certify ...\" -> justified, drop the now-pointless `_GoBack`\n# bookmark, and collapse the bookmark-split runs into a single run by\n# deleting and re-inserting the paragraph's full text in place. Note:\n# deleting the paragraph's whole range resets its alignment, so the\n# merge has to happen before the alignment is (re)applied.\nif ($certifyPara) {\n    if ($d.Bookmarks.Exists(\"_GoBack\")) {\n        $d.Bookmarks.Item(\"_GoBack\").Delete()\n    }\n\n    $r = $certifyPara.Range\n    $fullText = $r.Text\n    $r.Delete()\n    $r.InsertBefore($fullText)\n\n    $certifyPara.Alignment = $wdAlignParagraphJustify\n}\n"}
